$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 50.68470204858703

for ($row = 2; $row -le 7; $row++) {
    $ws.Range("N$row").Value = $newValue
}
